$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the time-slack (CPM) figures for the existing tasks (rows 2-18) ---
# Row 2 (task A)
$ws.Cells.Item(2,2).Value = 0
$ws.Cells.Item(2,3).Value = 12
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = 12
$ws.Cells.Item(2,6).Value = 0

# Row 3 (task B)
$ws.Cells.Item(3,2).Value = 12
$ws.Cells.Item(3,3).Value = 31
$ws.Cells.Item(3,4).Value = 12
$ws.Cells.Item(3,5).Value = 31
$ws.Cells.Item(3,6).Value = 0

# Row 4 (task C)
$ws.Cells.Item(4,2).Value = 31
$ws.Cells.Item(4,3).Value = 32
$ws.Cells.Item(4,4).Value = 31
$ws.Cells.Item(4,5).Value = 32
$ws.Cells.Item(4,6).Value = 0

# Row 5 (task D)
$ws.Cells.Item(5,2).Value = 32
$ws.Cells.Item(5,3).Value = 42
$ws.Cells.Item(5,4).Value = 32
$ws.Cells.Item(5,5).Value = 42
$ws.Cells.Item(5,6).Value = 0

# Row 6 (task E)
$ws.Cells.Item(6,2).Value = 42
$ws.Cells.Item(6,3).Value = 43
$ws.Cells.Item(6,4).Value = 42
$ws.Cells.Item(6,5).Value = 43
$ws.Cells.Item(6,6).Value = 0

# Row 7 (task F)
$ws.Cells.Item(7,2).Value = 43
$ws.Cells.Item(7,3).Value = 58
$ws.Cells.Item(7,4).Value = 43
$ws.Cells.Item(7,5).Value = 58
$ws.Cells.Item(7,6).Value = 0

# Row 8 (task G)
$ws.Cells.Item(8,2).Value = 12
$ws.Cells.Item(8,3).Value = 16
$ws.Cells.Item(8,4).Value = 53
$ws.Cells.Item(8,5).Value = 57
$ws.Cells.Item(8,6).Value = 41

# Row 9 (task H)
$ws.Cells.Item(9,2).Value = 16
$ws.Cells.Item(9,3).Value = 17
$ws.Cells.Item(9,4).Value = 53
$ws.Cells.Item(9,5).Value = 58
$ws.Cells.Item(9,6).Value = 42

# Row 10 (task I)
$ws.Cells.Item(10,2).Value = 58
$ws.Cells.Item(10,3).Value = 60
$ws.Cells.Item(10,4).Value = 58
$ws.Cells.Item(10,5).Value = 60
$ws.Cells.Item(10,6).Value = 0

# Row 11 (task J)
$ws.Cells.Item(11,2).Value = 58
$ws.Cells.Item(11,3).Value = 60
$ws.Cells.Item(11,4).Value = 60
$ws.Cells.Item(11,5).Value = 62
$ws.Cells.Item(11,6).Value = 2

# Row 12 (task K)
$ws.Cells.Item(12,2).Value = 58
$ws.Cells.Item(12,3).Value = 60
$ws.Cells.Item(12,4).Value = 60
$ws.Cells.Item(12,5).Value = 62
$ws.Cells.Item(12,6).Value = 2

# Row 13 (task L)
$ws.Cells.Item(13,2).Value = 60
$ws.Cells.Item(13,3).Value = 62
$ws.Cells.Item(13,4).Value = 62
$ws.Cells.Item(13,5).Value = 64
$ws.Cells.Item(13,6).Value = 2

# Row 14 (task N)
$ws.Cells.Item(14,2).Value = 60
$ws.Cells.Item(14,3).Value = 62
$ws.Cells.Item(14,4).Value = 60
$ws.Cells.Item(14,5).Value = 62
$ws.Cells.Item(14,6).Value = 0

# Row 15 (task M)
$ws.Cells.Item(15,2).Value = 60
$ws.Cells.Item(15,3).Value = 62
$ws.Cells.Item(15,4).Value = 62
$ws.Cells.Item(15,5).Value = 64
$ws.Cells.Item(15,6).Value = 2

# Row 16 (task O)
$ws.Cells.Item(16,2).Value = 62
$ws.Cells.Item(16,3).Value = 64
$ws.Cells.Item(16,4).Value = 62
$ws.Cells.Item(16,5).Value = 64
$ws.Cells.Item(16,6).Value = 0

# Row 17 (task P)
$ws.Cells.Item(17,2).Value = 64
$ws.Cells.Item(17,3).Value = 70
$ws.Cells.Item(17,4).Value = 64
$ws.Cells.Item(17,5).Value = 70
$ws.Cells.Item(17,6).Value = 0

# Row 18 (task Q)
$ws.Cells.Item(18,2).Value = 70
$ws.Cells.Item(18,3).Value = 73
$ws.Cells.Item(18,4).Value = 70
$ws.Cells.Item(18,5).Value = 73
$ws.Cells.Item(18,6).Value = 0

# --- Add the new task row (row 19, task R) ---
# Copy formatting (style + row height) from the row above, then overwrite values.
$ws.Range("A18:F18").Copy($ws.Range("A19:F19"))
$ws.Rows.Item(19).RowHeight = $ws.Rows.Item(18).RowHeight

$ws.Cells.Item(19,1).Value = "R"
$ws.Cells.Item(19,2).Value = 73
$ws.Cells.Item(19,3).Value = 78
$ws.Cells.Item(19,4).Value = 73
$ws.Cells.Item(19,5).Value = 78
$ws.Cells.Item(19,6).Value = 0

# --- Update the active selection to match the edited workbook ---
$ws.Range("H11").Select()
